$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Scenarios": flip the Execute flag from N to Y for rows 4-12
# ------------------------------------------------------------------
$scenarios = $wb.Worksheets.Item("Scenarios")
$scenarios.Range("A4:A12").Value2 = "Y"
$scenarios.Range("A4:A12").Select()

# ------------------------------------------------------------------
# Sheet "Parameters": various fixes
# ------------------------------------------------------------------
$params = $wb.Worksheets.Item("Parameters")

# Label the environment row
$params.Range("A8").Value2 = "ENV"

# retryLimit value should be stored as text "1" instead of a number
$params.Range("B12").NumberFormat = "@"
$params.Range("B12").Value2 = "1"

# Remove the old whole-number validation on the retry limit cell
$params.Range("B12").Validation.Delete()

# Apply a text number format to the cells that previously used the
# default/general format (this is what "fixes" the empty-row look)
$params.Range("A3:A6").NumberFormat = "@"
$params.Range("B3:B6").NumberFormat = "@"
$params.Range("A9:A11").NumberFormat = "@"
$params.Range("B9:B11").NumberFormat = "@"
$params.Range("A7:B7").NumberFormat = "@"
$params.Range("A8:B8").NumberFormat = "@"
$params.Columns.Item(2).NumberFormat = "@"

# Extend the text formatting down through row 45 so newly added rows
# pick up consistent formatting instead of being blank/unformatted
$params.Range("A13:A45").NumberFormat = "@"

$params.Range("B8").Select()

# ------------------------------------------------------------------
# Re-order / recolor the conditional formatting rules on Parameters
# ------------------------------------------------------------------
$b8 = $params.Range("B8")
$cf = $b8.FormatConditions
for ($i = $cf.Count; $i -ge 1; $i--) {
    $rule = $cf.Item($i)
    if ($rule.Formula1 -eq '=NOT(ISERROR(SEARCH("QA",B8)))') {
        $rule.Delete()
    }
}

$prodRule = $null
$cf2 = $b8.FormatConditions
for ($i = 1; $i -le $cf2.Count; $i++) {
    $rule = $cf2.Item($i)
    if ($rule.Formula1 -eq '=NOT(ISERROR(SEARCH("PROD",B8)))') {
        $prodRule = $rule
    }
}

$qaRule = $cf2.Add(9, 0, $null, $null, "QA", 0)
$qaRule.Font.Bold = $true
$qaRule.Interior.ThemeColor = 6

$wb.Save()
